$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("D2", "59.546.35"),
    @("E2", "  -1.03%  "),
    @("D3", "2.644.64"),
    @("E3", "  +1.21%  "),
    @("E4", "  +0.13%  "),
    @("D5", "'517.90"),
    @("E5", "  -0.56%  "),
    @("D6", "'146.99"),
    @("E6", "  -1.32%  "),
    @("D7", "'0.995"),
    @("E7", "  -0.34%  "),
    @("D8", "'0.577"),
    @("E8", "  +1.05%  "),
    @("D9", "2.668.82"),
    @("E9", "  +2.01%  "),
    @("D10", "'6.47"),
    @("E10", "  +3.22%  "),
    @("E11", "  +2.03%  "),
    @("D12", "'0.341"),
    @("E12", "  -0.38%  "),
    @("E13", "  -1.41%  "),
    @("D14", "3.109.21"),
    @("E14", "  +1.41%  "),
    @("D15", "59.508.49"),
    @("E15", "  -1.03%  "),
    @("D16", "'21.28"),
    @("E16", "  +0.33%  "),
    @("E17", "  +0.70%  "),
    @("D18", "2.667.58"),
    @("E18", "  +2.00%  "),
    @("D19", "'4.62"),
    @("E19", "  -0.23%  "),
    @("D20", "'346.01"),
    @("E20", "  +1.05%  "),
    @("E21", "  +1.65%  "),
    @("E22", "  +1.57%  "),
    @("D23", "'1.00"),
    @("E23", "  +0.22%  "),
    @("D24", "'61.08"),
    @("E24", "  +0.98%  "),
    @("D25", "'0.425"),
    @("E25", "  +1.14%  "),
    @("D26", "2.771.34"),
    @("E26", "  +1.67%  "),
    @("D27", "'0.163"),
    @("E27", "  +0.78%  "),
    @("D28", "'0.993"),
    @("D29", "0.0₃0822"),
    @("E29", "  +1.79%  "),
    @("D30", "'7.20"),
    @("E30", "  +1.74%  "),
    @("D31", "'0.997"),
    @("E31", "  -0.31%  "),
    @("E32", "  +8.52%  "),
    @("D33", "'19.08"),
    @("E33", "  +0.65%  "),
    @("E34", "  -0.13%  "),
    @("B35", "Monero"),
    @("C35", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"),
    @("D35", "'149.97"),
    @("E35", "  +0.33%  "),
    @("B36", "SuiNetwork"),
    @("C36", "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"),
    @("D36", "'1.06"),
    @("E36", "  +15.92%  "),
    @("E37", "  +2.95%  "),
    @("E38", "  +2.86%  "),
    @("D39", "'0.871"),
    @("E39", "  +0.95%  "),
    @("D40", "'36.55"),
    @("E40", "  +0.28%  "),
    @("E41", "  +3.50%  "),
    @("E42", "  +0.15%  "),
    @("D43", "'285.67"),
    @("E43", "  -1.09%  "),
    @("D44", "'0.621"),
    @("E44", "  -0.27%  "),
    @("D45", "'0.0998"),
    @("E45", "  -0.60%  "),
    @("D46", "'0.993"),
    @("E46", "  -0.54%  "),
    @("D47", "'19.77"),
    @("E47", "  +1.59%  "),
    @("E48", "  -0.23%  "),
    @("E49", "  +1.41%  "),
    @("E50", "  +1.53%  "),
    @("B51", "WhiteBITCoin"),
    @("C51", "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"),
    @("D51", "'10.26"),
    @("E51", "  -1.27%  ")
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}
